# Arrumei a classificação de alguns tweets
$wb = $excel.ActiveWorkbook
$wsTr = $wb.Worksheets.Item("Treinamento")
$wsT  = $wb.Worksheets.Item("Teste")

# --- Reclassify some tweets in "Treinamento" (B column 1 -> 0) ---
$fixRows = @(7, 9, 22, 43, 47, 55, 102, 144, 168, 180, 221, 224, 239, 240, 250)
foreach ($r in $fixRows) {
    $wsTr.Cells.Item($r, 2).Value = 0
}

# --- Underline the tweet text for rows that got reclassified as highlighted (A7, A12) ---
$wsTr.Range("A7").Font.Underline = 2
$wsTr.Range("A12").Font.Underline = 2

# --- Add a new (currently empty) row marker at the bottom of "Treinamento", styled like A7/A12 ---
$wsTr.Range("A303").Font.Underline = 2

# --- Unify the header style of "Teste"!A1 with "Treinamento"!A1 (bold + underline + border) ---
$wsT.Range("A1").Font.Underline = 2

# --- Sheet view / zoom / selection updates ---
$wsT.Select()
$excel.ActiveWindow.Zoom = 70
$wsT.Range("A14").Select()

$wsTr.Select()
$excel.ActiveWindow.Zoom = 70
$wsTr.Range("A4").Select()
